$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the season-record columns: AD=Wins, AE=Losses, AF=Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Give the new headers the same style as the existing header row (bold, bordered, centered)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record (70 wins, 92 losses, 0 ties) for every player row
for ($row = 2; $row -le 54; $row++) {
    $ws.Cells.Item($row, 30).Value = 70
    $ws.Cells.Item($row, 31).Value = 92
    $ws.Cells.Item($row, 32).Value = 0
}
